$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (push existing data rows 2-21 down to 3-22);
# Excel carries the formatting of the row above into the freshly inserted row.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the "Column foundation" record that used
# to be folded into the header (it becomes id 1; everything below it renumbers).
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Column foundation"
$ws.Range("C2").Value = "piling and foundations"
$ws.Range("D2").Value = 60

# Renumber the id column for all the rows that shifted down (old id N -> row N+2)
for ($r = 3; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# The header cell in B1 changes from "Column foundation" to the column title "Assembly",
# and loses the leftover (no-op) formatting override it used to carry.
$ws.Range("B1").Value = "Assembly"
$ws.Range("B1").ClearFormats()

# Move the visible selection to reflect where the author was last working
$ws.Range("C10").Select() | Out-Null

$wb.Save() | Out-Null
